# Reorders the D, J, K, L, M, O, P values across rows 2-14 of the active sheet.
# The row identity (A, B, C, E, F, G, H, I, N, Q, R) stays constant per row;
# only the "record" formed by columns D, J, K, L, M, O, P moves between rows,
# per the mapping below (new row -> source/original row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new row number -> original row number that supplies its D/J/K/L/M/O/P values
$map = @{
    2  = 5
    3  = 4
    4  = 13
    5  = 8
    6  = 3
    7  = 14
    8  = 12
    9  = 10
    10 = 11
    11 = 9
    12 = 7
    13 = 2
    14 = 6
}

# Snapshot original values for the columns that move, keyed by original row.
$cols = @("D", "J", "K", "L", "M", "O", "P")
$snapshot = @{}
foreach ($r in 2..14) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Apply the permutation: each new row gets the snapshot values of its mapped source row.
foreach ($newRow in 2..14) {
    $srcRow = $map[$newRow]
    $rowVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $rowVals[$c]
    }
}
